$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at C for "Matching Filter"
$ws.Columns("C:C").Insert()

# 2. Insert a new row at 2 for "Les Acacias" (pushes the Saint Avertin row down to 3)
$ws.Rows("2:2").Insert()

# 3. Insert a new row at 4 for "L'Islette - Fondettes" (between Saint Avertin row 3
#    and the Mignardiere row, which becomes row 5)
$ws.Rows("4:4").Insert()

# New header for the inserted column
$ws.Range("C1").Value = 'Matching Filter'

# Row/cell values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 4.62
$ws.Range("C2").Value = 'tourism=camp_site'
$ws.Range("D2").Value = 'Les Acacias'
$ws.Range("E2").Value = 'https://www.camping-tours.fr/'
$ws.Range("F2").Value = '+33 2 47 44 08 16'
$ws.Range("H2").Value = '{''addr:city'': ''La Ville-aux-Dames'', ''addr:postcode'': ''37700'', ''addr:street'': ''Rue Berthe Morisot'', ''barrier'': ''fence'', ''caravans'': ''yes'', ''email'': ''contact@camplvad.com'', ''internet_access'': ''yes'', ''internet_access:fee'': ''no'', ''name'': ''Les Acacias'', ''phone'': ''+33 2 47 44 08 16'', ''stars'': ''3'', ''tents'': ''yes'', ''tourism'': ''camp_site'', ''website'': ''https://www.camping-tours.fr/''}'
$ws.Range("I2").Value = 47.4020858
$ws.Range("J2").Value = 0.7801299
$ws.Range("A3").Value = 1.57
$ws.Range("B3").Value = 1.57
$ws.Range("C3").Value = 'tourism=camp_site'
$ws.Range("D3").Value = 'Aire camping-car de Saint Avertin'
$ws.Range("E3").Value = 'https://www.onlypark.fr/aire-camping-car-de-st-avertin/'
$ws.Range("F3").Value = '+33 2 47 27 87 47'
$ws.Range("H3").Value = '{''addr:city'': ''Saint-Avertin'', ''capacity'': ''20'', ''caravans'': ''yes'', ''charge'': ''19 EUR'', ''charge:conditional'': ''12 EUR @ (Sep-Jun)'', ''contact:email'': ''campingtoursvaldeloire@onlycamp.fr'', ''contact:phone'': ''+33 2 47 27 87 47'', ''contact:website'': ''https://www.onlypark.fr/aire-camping-car-de-st-avertin/'', ''drinking_water'': ''yes'', ''fee'': ''yes'', ''internet_access'': ''wlan'', ''name'': ''Aire camping-car de Saint Avertin'', ''network'': ''Onlypark'', ''operator'': ''onlycamp'', ''power_supply'': ''yes'', ''power_supply:charge'': ''3 EUR/4 hours'', ''power_supply:fee'': ''yes'', ''power_supply:maxcurrent'': ''10'', ''sanitary_dump_station'': ''yes'', ''sanitary_dump_station:charge'': ''3 EUR/20 minutes'', ''sanitary_dump_station:fee'': ''yes'', ''shower'': ''yes'', ''stars'': ''4'', ''tents'': ''yes'', ''toilets'': ''no'', ''tourism'': ''camp_site'', ''water_point'': ''yes'', ''wheelchair'': ''yes''}'
$ws.Range("I3").Value = 47.3708862
$ws.Range("J3").Value = 0.7243202
$ws.Range("A4").Value = 9.91
$ws.Range("B4").Value = 4.44
$ws.Range("C4").Value = 'tourism=camp_site'
$ws.Range("D4").Value = 'L''Islette - Fondettes'
$ws.Range("H4").Value = '{''name'': "L''Islette - Fondettes", ''tourism'': ''camp_site''}'
$ws.Range("I4").Value = 47.3892756
$ws.Range("J4").Value = 0.5959279
$ws.Range("A5").Value = 11.11
$ws.Range("B5").Value = 0.08
$ws.Range("C5").Value = 'tourism=camp_site'
$ws.Range("D5").Value = 'Camping La Mignardière'
$ws.Range("E5").Value = 'https://www.mignardiere.com/'
$ws.Range("F5").Value = '+33 2 47 73 31 00'
$ws.Range("H5").Value = '{''addr:city'': ''Ballan-Miré'', ''addr:housenumber'': ''22'', ''addr:postcode'': ''37510'', ''addr:street'': ''Avenue des Aubépines'', ''cabins'': ''yes'', ''capacity:caravans'': ''114'', ''capacity:tents'': ''114'', ''caravans'': ''yes'', ''drinking_water'': ''yes'', ''motorhome'': ''yes'', ''name'': ''Camping La Mignardière'', ''phone'': ''+33 2 47 73 31 00'', ''sanitary_dump_station'': ''yes'', ''shower'': ''yes'', ''stars'': ''4'', ''tents'': ''yes'', ''toilets'': ''yes'', ''tourism'': ''camp_site'', ''washing_machine'': ''yes'', ''website'': ''https://www.mignardiere.com/''}'
$ws.Range("I5").Value = 47.3557614
$ws.Range("J5").Value = 0.6332265


# Clear all cell-level formatting (bold header font + borders + centered alignment)
# so every cell reverts to the default (unstyled) cell format.
$ws.Cells.ClearFormats()
